$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 301, pushing existing rows 301-311 down to 302-312
$ws.Rows("301:301").Insert()

# Populate the newly inserted row 301 with the new data record
$ws.Range("A301").Value = 9
$ws.Range("B301").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C301").Value = "Metropolitana"
$ws.Range("D301").Value = 44509
$ws.Range("E301").Value = 13
$ws.Range("F301").Value = 100112024
$ws.Range("G301").Value = "Choclo"
$ws.Range("H301").Value = "Choclero"
$ws.Range("I301").Value = "Primera"
$ws.Range("J301").Value = 34
$ws.Range("K301").Value = 25000
$ws.Range("L301").Value = 26000
$ws.Range("M301").Value = 25500
$ws.Range("N301").Value = "$/caja 50 unidades"
$ws.Range("O301").Value = "Argentina"
$ws.Range("P301").Value = 510
$ws.Range("Q301").Value = 50
$ws.Range("R301").Value = "Hortaliza"
